# Update registration counters (Inscritos, Pagos, Inscrições homologadas)
# in the "IFMG - Processo Seletivo Cursos Técnicos Integrados 2025/1" sheet
# to reflect the latest count as of the commit date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 85
$ws.Range("H5").Value = 85
$ws.Range("E10").Value = 452
$ws.Range("F10").Value = 228
$ws.Range("H10").Value = 228
$ws.Range("E11").Value = 313
$ws.Range("F11").Value = 175
$ws.Range("H11").Value = 175
$ws.Range("E12").Value = 444
$ws.Range("F12").Value = 245
$ws.Range("H12").Value = 245
$ws.Range("F13").Value = 63
$ws.Range("H13").Value = 63
$ws.Range("F14").Value = 59
$ws.Range("H14").Value = 59
$ws.Range("F15").Value = 61
$ws.Range("H15").Value = 61
$ws.Range("E16").Value = 186
$ws.Range("F16").Value = 95
$ws.Range("H16").Value = 95
$ws.Range("F17").Value = 45
$ws.Range("H17").Value = 45
$ws.Range("F21").Value = 71
$ws.Range("H21").Value = 71
$ws.Range("E22").Value = 158
$ws.Range("F22").Value = 80
$ws.Range("H22").Value = 80
$ws.Range("F23").Value = 82
$ws.Range("H23").Value = 82
$ws.Range("F24").Value = 105
$ws.Range("H24").Value = 105
$ws.Range("E25").Value = 242
$ws.Range("F25").Value = 116
$ws.Range("H25").Value = 116
$ws.Range("F26").Value = 87
$ws.Range("H26").Value = 87
$ws.Range("F27").Value = 144
$ws.Range("H27").Value = 144
$ws.Range("F28").Value = 67
$ws.Range("H28").Value = 67
$ws.Range("E29").Value = 156
$ws.Range("F29").Value = 88
$ws.Range("H29").Value = 88
$ws.Range("E30").Value = 193
$ws.Range("F30").Value = 113
$ws.Range("H30").Value = 113
$ws.Range("E32").Value = 173
$ws.Range("F32").Value = 101
$ws.Range("H32").Value = 101
$ws.Range("E33").Value = 262
$ws.Range("F33").Value = 136
$ws.Range("H33").Value = 136
$ws.Range("F34").Value = 128
$ws.Range("H34").Value = 128
$ws.Range("F35").Value = 85
$ws.Range("H35").Value = 85
$ws.Range("F37").Value = 72
$ws.Range("H37").Value = 72
$ws.Range("E40").Value = 240
$ws.Range("F40").Value = 111
$ws.Range("H40").Value = 111
$ws.Range("E41").Value = 364
$ws.Range("F41").Value = 172
$ws.Range("H41").Value = 172
$ws.Range("E42").Value = 332
$ws.Range("F42").Value = 183
$ws.Range("H42").Value = 183
$ws.Range("F44").Value = 146
$ws.Range("H44").Value = 146
$ws.Range("E45").Value = 129
$ws.Range("F45").Value = 64
$ws.Range("H45").Value = 64
$ws.Range("E46").Value = 286
$ws.Range("F46").Value = 158
$ws.Range("H46").Value = 158
$ws.Range("E47").Value = 409
$ws.Range("F47").Value = 201
$ws.Range("H47").Value = 201
$ws.Range("F48").Value = 78
$ws.Range("H48").Value = 78
$ws.Range("E49").Value = 267
$ws.Range("F49").Value = 115
$ws.Range("H49").Value = 115
$ws.Range("F50").Value = 104
$ws.Range("H50").Value = 104
$ws.Range("F51").Value = 90
$ws.Range("H51").Value = 90

$wb.Save()
